# Atualização de bases das ligas, do dia: 23-02-2024 às 23:34
#
# Rows 164 and 165 swap their data (columns B:AC), and rows 175-179
# cyclically rotate their data (columns B:AC): row 179's data moves to
# row 175, and rows 175,176,177,178 each shift down into 176,177,178,179.
# Column A (the row index) is left untouched in both groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

function Get-RowValues($row, $c1, $c2) {
    $vals = @()
    for ($c = $c1; $c -le $c2; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $c1, $c2, $vals) {
    $i = 0
    for ($c = $c1; $c -le $c2; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

# --- Swap rows 164 and 165 (columns B:AC) ---
$row164 = Get-RowValues 164 $firstCol $lastCol
$row165 = Get-RowValues 165 $firstCol $lastCol

Set-RowValues 164 $firstCol $lastCol $row165
Set-RowValues 165 $firstCol $lastCol $row164

# --- Cyclic shift rows 175..179 (columns B:AC) ---
# new(175) = old(179); new(176) = old(175); new(177) = old(176);
# new(178) = old(177); new(179) = old(178)
$row175 = Get-RowValues 175 $firstCol $lastCol
$row176 = Get-RowValues 176 $firstCol $lastCol
$row177 = Get-RowValues 177 $firstCol $lastCol
$row178 = Get-RowValues 178 $firstCol $lastCol
$row179 = Get-RowValues 179 $firstCol $lastCol

Set-RowValues 175 $firstCol $lastCol $row179
Set-RowValues 176 $firstCol $lastCol $row175
Set-RowValues 177 $firstCol $lastCol $row176
Set-RowValues 178 $firstCol $lastCol $row177
Set-RowValues 179 $firstCol $lastCol $row178
